$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to be stored as text, matching the
# original inline-string cell type, so values like "42.00" or "235.84"
# are not silently coerced into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "95.016.75"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "3.602.91"
$ws.Range("E3").Value = "  +4.73%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "235.84"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Value = "656.38"
$ws.Range("E6").Value = "  +4.75%  "
$ws.Range("D7").Value = "1.45"
$ws.Range("E7").Value = "  +1.29%  "
$ws.Range("D8").Value = "0.399"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "0.989"
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D11").Value = "3.599.89"
$ws.Range("E11").Value = "  +4.67%  "
$ws.Range("D12").Value = "0.201"
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").Value = "42.00"
$ws.Range("E13").Value = "  -3.29%  "
$ws.Range("D14").Value = "6.32"
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("D15").Value = "4.291.54"
$ws.Range("E15").Value = "  +4.56%  "
$ws.Range("D16").Value = "94.894.84"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("D17").Value = "0.0000251"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").Value = "3.592.39"
$ws.Range("E18").Value = "  +4.23%  "
$ws.Range("D19").Value = "7.91"
$ws.Range("E19").Value = "  -4.51%  "
$ws.Range("D20").Value = "12.76"
$ws.Range("E20").Value = "  +8.74%  "
$ws.Range("D21").Value = "17.87"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("D22").Value = "3.54"
$ws.Range("E22").Value = "  +4.49%  "
$ws.Range("D23").Value = "0.479"
$ws.Range("E23").Value = "  -2.69%  "
$ws.Range("D24").Value = "501.29"
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").Value = "0.0000196"
$ws.Range("E25").Value = "  +6.51%  "
$ws.Range("D26").Value = "6.56"
$ws.Range("E26").Value = "  -3.17%  "
$ws.Range("D27").Value = "94.75"
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("D28").Value = "3.798.17"
$ws.Range("E28").Value = "  +4.74%  "
$ws.Range("D29").Value = "12.40"
$ws.Range("E29").Value = "  +2.29%  "
$ws.Range("D30").Value = "3.09"
$ws.Range("E30").Value = "  +8.74%  "
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "11.17"
$ws.Range("E32").Value = "  -1.54%  "
$ws.Range("D33").Value = "0.138"
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("D34").Value = "0.995"
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("D35").Value = "32.03"
$ws.Range("E35").Value = "  +8.93%  "
$ws.Range("D36").Value = "0.175"
$ws.Range("E36").Value = "  -2.08%  "
$ws.Range("D37").Value = "0.554"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "563.72"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").Value = "8.00"
$ws.Range("E39").Value = "  +6.18%  "
$ws.Range("D40").Value = "1.46"
$ws.Range("E40").Value = "  +2.20%  "
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "0.149"
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("D43").Value = "0.907"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("D44").Value = "34.98"
$ws.Range("E44").Value = "  +43.24%  "
$ws.Range("D45").Value = "1.73"
$ws.Range("E45").Value = "  +2.01%  "
$ws.Range("D46").Value = "23.67"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").Value = "5.57"
$ws.Range("E47").Value = "  +1.03%  "
$ws.Range("D48").Value = "2.22"
$ws.Range("E48").Value = "  +4.83%  "
$ws.Range("D49").Value = "0.0410"
$ws.Range("E49").Value = "  -3.12%  "
$ws.Range("D50").Value = "3.53"
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("D51").Value = "53.39"
$ws.Range("E51").Value = "  +0.42%  "
